# Generate Report for Handoff
# Updates the localization-status workbook: the e27be0ea file is now
# "Ready for handoff" (previously "Handed back: in sync with en-US"),
# and the zh-cn / de-de "Latest Handoff Datetime" values are refreshed
# to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"

# Overview sheet: row 3 is the e27be0ea-aa8b-4bbb-a1e7-5f6063e618fd.md entry
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff

# zh-cn sheet: row 3 Status -> Ready for handoff, row 2 & 3 Latest Handoff Datetime updated
$wsZhCn.Range("B3").Value = $readyForHandoff
$wsZhCn.Range("D2").Value = "2016-03-04 04:10:30"
$wsZhCn.Range("D3").Value = "2016-03-04 04:10:30"

# de-de sheet: row 3 Status -> Ready for handoff, row 2 & 3 Latest Handoff Datetime updated
$wsDeDe.Range("B3").Value = $readyForHandoff
$wsDeDe.Range("D2").Value = "2016-03-04 04:10:46"
$wsDeDe.Range("D3").Value = "2016-03-04 04:10:46"
